$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Reroll Packs")
$ws2 = $wb.Worksheets.Item("Blad1")
$ws3 = $wb.Worksheets.Item("Reroll History")

# --- Reroll History ("Pocket" tab): add new promo pull rows (10-16) ---
# Populate data rows first so new card-name / type shared strings are
# allocated in the same order the original author typed them in, before
# the new "Code" header (kept last) is written.

$ws3.Range("A57").Value = 10
$ws3.Range("B57").Value = 'Exeggutor ex'
$ws3.Range("C57").Value = 'ex'
$ws3.Range("D57").Value = 1
$ws3.Range("E57").Formula = '=IF(D57-1 = 0, "", D57-1)'

$ws3.Range("A58").Value = 10
$ws3.Range("B58").Value = 'Blastoise ex'
$ws3.Range("C58").Value = 'ex'
$ws3.Range("D58").Value = 1
$ws3.Range("E58").Formula = '=IF(D58-1 = 0, "", D58-1)'

$ws3.Range("A59").Value = 10
$ws3.Range("B59").Value = 'Gengar ex'
$ws3.Range("C59").Value = 'ex'
$ws3.Range("D59").Value = 1
$ws3.Range("E59").Formula = '=IF(D59-1 = 0, "", D59-1)'

$ws3.Range("A60").Value = 10
$ws3.Range("B60").Value = 'Wigglytuff ex'
$ws3.Range("C60").Value = 'ex'
$ws3.Range("D60").Value = 1
$ws3.Range("E60").Formula = '=IF(D60-1 = 0, "", D60-1)'

$ws3.Range("A61").Value = 10
$ws3.Range("B61").Value = 'Charmander'
$ws3.Range("C61").Value = 'Special_Art_Pokémon'
$ws3.Range("D61").Value = 1
$ws3.Range("E61").Formula = '=IF(D61-1 = 0, "", D61-1)'

$ws3.Range("A62").Value = 10
$ws3.Range("B62").Value = 'Gyarados'
$ws3.Range("C62").Value = 'Special_Art_Pokémon'
$ws3.Range("D62").Value = 1
$ws3.Range("E62").Formula = '=IF(D62-1 = 0, "", D62-1)'

$ws3.Range("A63").Value = 10
$ws3.Range("B63").Value = 'Diglett'
$ws3.Range("C63").Value = 'Special_Art_Pokémon'
$ws3.Range("D63").Value = 1
$ws3.Range("E63").Formula = '=IF(D63-1 = 0, "", D63-1)'

$ws3.Range("A64").Value = 11
$ws3.Range("B64").Value = 'Exeggutor ex'
$ws3.Range("C64").Value = 'ex'
$ws3.Range("D64").Value = 1
$ws3.Range("E64").Formula = '=IF(D64-1 = 0, "", D64-1)'

$ws3.Range("A65").Value = 11
$ws3.Range("B65").Value = 'Zapdos ex'
$ws3.Range("C65").Value = 'ex'
$ws3.Range("D65").Value = 1
$ws3.Range("E65").Formula = '=IF(D65-1 = 0, "", D65-1)'

$ws3.Range("A66").Value = 11
$ws3.Range("B66").Value = 'Charmander'
$ws3.Range("C66").Value = 'Special_Art_Pokémon'
$ws3.Range("D66").Value = 1
$ws3.Range("E66").Formula = '=IF(D66-1 = 0, "", D66-1)'

$ws3.Range("A67").Value = 11
$ws3.Range("B67").Value = 'Gyarados'
$ws3.Range("C67").Value = 'Special_Art_Pokémon'
$ws3.Range("D67").Value = 1
$ws3.Range("E67").Formula = '=IF(D67-1 = 0, "", D67-1)'

$ws3.Range("A68").Value = 12
$ws3.Range("B68").Value = 'Exeggutor ex'
$ws3.Range("C68").Value = 'ex'
$ws3.Range("D68").Value = 1
$ws3.Range("E68").Formula = '=IF(D68-1 = 0, "", D68-1)'

$ws3.Range("A69").Value = 12
$ws3.Range("B69").Value = 'Charmander'
$ws3.Range("C69").Value = 'Special_Art_Pokémon'
$ws3.Range("D69").Value = 1
$ws3.Range("E69").Formula = '=IF(D69-1 = 0, "", D69-1)'

$ws3.Range("A70").Value = 12
$ws3.Range("B70").Value = 'Gyarados'
$ws3.Range("C70").Value = 'Special_Art_Pokémon'
$ws3.Range("D70").Value = 1
$ws3.Range("E70").Formula = '=IF(D70-1 = 0, "", D70-1)'

$ws3.Range("A71").Value = 12
$ws3.Range("B71").Value = 'Snorlax'
$ws3.Range("C71").Value = 'Special_Art_Pokémon'
$ws3.Range("D71").Value = 1
$ws3.Range("E71").Formula = '=IF(D71-1 = 0, "", D71-1)'

$ws3.Range("A72").Value = 13
$ws3.Range("B72").Value = 'Exeggutor ex'
$ws3.Range("C72").Value = 'ex'
$ws3.Range("D72").Value = 1
$ws3.Range("E72").Formula = '=IF(D72-1 = 0, "", D72-1)'

$ws3.Range("A73").Value = 13
$ws3.Range("B73").Value = 'Charmander'
$ws3.Range("C73").Value = 'Special_Art_Pokémon'
$ws3.Range("D73").Value = 1
$ws3.Range("E73").Formula = '=IF(D73-1 = 0, "", D73-1)'

$ws3.Range("A74").Value = 13
$ws3.Range("B74").Value = 'Diglett'
$ws3.Range("C74").Value = 'Special_Art_Pokémon'
$ws3.Range("D74").Value = 1
$ws3.Range("E74").Formula = '=IF(D74-1 = 0, "", D74-1)'

$ws3.Range("A75").Value = 14
$ws3.Range("B75").Value = 'Exeggutor ex'
$ws3.Range("C75").Value = 'ex'
$ws3.Range("D75").Value = 1
$ws3.Range("E75").Formula = '=IF(D75-1 = 0, "", D75-1)'

$ws3.Range("A76").Value = 14
$ws3.Range("B76").Value = 'Wigglytuff ex'
$ws3.Range("C76").Value = 'ex'
$ws3.Range("D76").Value = 1
$ws3.Range("E76").Formula = '=IF(D76-1 = 0, "", D76-1)'

$ws3.Range("A77").Value = 14
$ws3.Range("B77").Value = 'Charmander'
$ws3.Range("C77").Value = 'Special_Art_Pokémon'
$ws3.Range("D77").Value = 1
$ws3.Range("E77").Formula = '=IF(D77-1 = 0, "", D77-1)'

$ws3.Range("A78").Value = 14
$ws3.Range("B78").Value = 'Cubone'
$ws3.Range("C78").Value = 'Special_Art_Pokémon'
$ws3.Range("D78").Value = 1
$ws3.Range("E78").Formula = '=IF(D78-1 = 0, "", D78-1)'

$ws3.Range("A79").Value = 14
$ws3.Range("B79").Value = 'Arcanine ex'
$ws3.Range("C79").Value = 'Full_Art_ex'
$ws3.Range("D79").Value = 1
$ws3.Range("E79").Formula = '=IF(D79-1 = 0, "", D79-1)'

$ws3.Range("A80").Value = 15
$ws3.Range("B80").Value = 'Exeggutor ex'
$ws3.Range("C80").Value = 'ex'
$ws3.Range("D80").Value = 1
$ws3.Range("E80").Formula = '=IF(D80-1 = 0, "", D80-1)'

$ws3.Range("A81").Value = 15
$ws3.Range("B81").Value = 'Marowak ex'
$ws3.Range("C81").Value = 'ex'
$ws3.Range("D81").Value = 1
$ws3.Range("E81").Formula = '=IF(D81-1 = 0, "", D81-1)'

$ws3.Range("A82").Value = 15
$ws3.Range("B82").Value = 'Charmander'
$ws3.Range("C82").Value = 'Special_Art_Pokémon'
$ws3.Range("D82").Value = 1
$ws3.Range("E82").Formula = '=IF(D82-1 = 0, "", D82-1)'

$ws3.Range("A83").Value = 15
$ws3.Range("B83").Value = 'Squirtle'
$ws3.Range("C83").Value = 'Special_Art_Pokémon'
$ws3.Range("D83").Value = 1
$ws3.Range("E83").Formula = '=IF(D83-1 = 0, "", D83-1)'

$ws3.Range("A84").Value = 15
$ws3.Range("B84").Value = 'Alakazam'
$ws3.Range("C84").Value = 'Special_Art_Pokémon'
$ws3.Range("D84").Value = 1
$ws3.Range("E84").Formula = '=IF(D84-1 = 0, "", D84-1)'

$ws3.Range("A85").Value = 15
$ws3.Range("B85").Value = 'Eevee'
$ws3.Range("C85").Value = 'Special_Art_Pokémon'
$ws3.Range("D85").Value = 1
$ws3.Range("E85").Formula = '=IF(D85-1 = 0, "", D85-1)'

$ws3.Range("A86").Value = 16
$ws3.Range("B86").Value = 'Exeggutor ex'
$ws3.Range("C86").Value = 'ex'
$ws3.Range("D86").Value = 1
$ws3.Range("E86").Formula = '=IF(D86-1 = 0, "", D86-1)'

$ws3.Range("A87").Value = 16
$ws3.Range("B87").Value = 'Starmie ex'
$ws3.Range("C87").Value = 'ex'
$ws3.Range("D87").Value = 1

$ws3.Range("A88").Value = 16
$ws3.Range("B88").Value = 'Charizard ex'
$ws3.Range("C88").Value = 'ex'
$ws3.Range("D88").Value = 1

$ws3.Range("A89").Value = 16
$ws3.Range("B89").Value = 'Charmander'
$ws3.Range("C89").Value = 'Special_Art_Pokémon'
$ws3.Range("D89").Value = 1

$ws3.Range("A90").Value = 16
$ws3.Range("B90").Value = 'Nidoqueen'
$ws3.Range("C90").Value = 'Special_Art_Pokémon'
$ws3.Range("D90").Value = 1

$ws3.Range("A91").Value = 16
$ws3.Range("B91").Value = 'Misty'
$ws3.Range("C91").Value = 'Full_Art_Trainer'
$ws3.Range("D91").Value = 1

$ws3.Range("A92").Value = 16
$ws3.Range("B92").Value = 'Arcanine ex'
$ws3.Range("C92").Value = 'Full_Art_ex'
$ws3.Range("D92").Value = 1

$ws3.Range("A93").Value = 16
$ws3.Range("B93").Value = 'Pikachu ex'
$ws3.Range("C93").Value = 'Full_Art_ex'
$ws3.Range("D93").Value = 1

$ws3.Range("A94").Value = 16
$ws3.Range("B94").Value = 'Pikachu ex'
$ws3.Range("C94").Value = 'Immersive_Art'
$ws3.Range("D94").Value = 1

# --- Reroll History: new probability columns (H/I) for rows 2-5, and the new "Code" header (F1) ---
$ws3.Range('H2').Value = 0.85
$ws3.Range('H3').Formula = '=3/20'
$ws3.Range('I3').Formula = '=H3 * (3/20)'
$ws3.Range('H4').Formula = '=1/4'
$ws3.Range('I4').Formula = '=H4 * (1/4)'
$ws3.Range('H5').Formula = '=1/16'
$ws3.Range('F1').Value = 'Code'

# --- View / selection state to match the saved workbook ---
$ws1.Activate()
$ws1.Range('A2').Select()
$ws1.Range('A15').Select()

$ws2.Activate()
$ws2.Range('B10').Select()

$ws3.Activate()
$ws3.Range('I4').Select()
